$d = $word.ActiveDocument

# --- Change 1: remove italics from the first "fridlysta" comment paragraph
#     and append a period to its text. ---
$rng1 = $d.Content
$oldText1 = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"
$newText1 = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen."
$rng1.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $newText1, 1)
$rng1.Italic = 0

# --- Change 2: trim the trailing space on the 6.4 paragraph text ---
$rng2 = $d.Content
$oldText2 = "Certifikatsinnehavaren ska skydda sällsynta arter och hotade arter samt deras livsmiljöer inom skogsbruksenheten. Det ska ske genom avsättningar, andra skyddade områden och genom att skapa konnektivitet och/eller genom andra direkta åtgärder som gynnar dessa arters överlevnad och livskraft. Åtgärderna ska stå i förhållande till brukandets skala, intensitet och risk, samt till sällsynta och hotade arters bevarandestatus och ekologiska krav. Certifikatsinnehavaren ska beakta den geografiska spridningen och ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas. "
$newText2 = "Certifikatsinnehavaren ska skydda sällsynta arter och hotade arter samt deras livsmiljöer inom skogsbruksenheten. Det ska ske genom avsättningar, andra skyddade områden och genom att skapa konnektivitet och/eller genom andra direkta åtgärder som gynnar dessa arters överlevnad och livskraft. Åtgärderna ska stå i förhållande till brukandets skala, intensitet och risk, samt till sällsynta och hotade arters bevarandestatus och ekologiska krav. Certifikatsinnehavaren ska beakta den geografiska spridningen och ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas."
$rng2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 1)

# --- Change 3: delete the "6.4.1 Följande biotoper..." paragraph and the
#     "b) nyckelbiotoper..." paragraph that follows it, then renumber the
#     remaining "6.4.1 " heading (before "Bevarandeåtgärder...") to "6.4.3 ". ---
$rng3 = $d.Content
$rng3.Find.Execute("Följande biotoper undantas från alla skogsbruksåtgärder", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$pFirst = $rng3.Paragraphs(1)

$rng4 = $d.Content
$rng4.Find.Execute("b) nyckelbiotoper enligt Skogsstyrelsens definition och metod (1995)", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$pLast = $rng4.Paragraphs(1)

$delRange = $d.Range($pFirst.Range.Start, $pLast.Range.End)
$delRange.Delete()

$rng5 = $d.Content
$rng5.Find.Execute("6.4.1 ", $true, $false, $false, $false, $false, $true, 1, $false, "6.4.3 ", 1)

# --- Change 4: update the date in the first-page header. ---
$sec = $d.Sections(1)
$hdr = $sec.Headers(2)
$hdr.Range.Find.Execute("2023-10-22", $false, $false, $false, $false, $false, $true, 1, $false, "2023-10-25", 1)
